# Progress on the implementation chapter of the dissertation:
# fill in the time-log rows that were left blank (rows 21-36) with the
# dates / start & end times / summary text for the new work that was logged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (date serial, start time, end time, summary)
$entries = @(
  @{ Row = 21; Date = 45239; Start = 0.55555555555555558; End = 0.7284722222222223;  Summary = "Work on shader preprocessing" },
  @{ Row = 22; Date = 45241; Start = 0.30555555555555552; End = 0.74097222222222225; Summary = "Completed shader preprocessing" },
  @{ Row = 23; Date = 45250; Start = 0.95833333333333337; End = 0.083333333333333329; Summary = "Multi-pass/multi-buffer api design & research" },
  @{ Row = 24; Date = 45257; Start = 0.91666666666666663; End = 0.10416666666666667; Summary = "Multi-pass/multi-buffer api implementation" },
  @{ Row = 25; Date = 45261; Start = 0.375;                End = 0.57222222222222219; Summary = "Multi-pass/multi-buffer api implementation & widget improvements" },
  @{ Row = 26; Date = 45262; Start = 0.32708333333333334; End = 0.5180555555555556;  Summary = "Multi-pass/multi-buffer api implementation" },
  @{ Row = 27; Date = 45265; Start = 0.16666666666666666; End = 0.41666666666666669; Summary = "Multi-pass rendering" },
  @{ Row = 28; Date = 45266; Start = 0.38541666666666669; End = 0.53472222222222221; Summary = "Renderdoc integration & geometry shaders" },
  @{ Row = 29; Date = 45267; Start = 0.5;                 End = 0.72916666666666663; Summary = "Finished implementing textures & updated documentation" },
  @{ Row = 30; Date = 45302; Start = 0.125;                End = 0.41666666666666669; Summary = "Fixed geometry shaders and various other bugs. Prepared a new release" },
  @{ Row = 31; Date = 45303; Start = 0.33333333333333331; End = 0.89583333333333337; Summary = "Camera bugfixes & UI library implementation" },
  @{ Row = 32; Date = 45312; Start = 0.10416666666666667; End = 0.54166666666666663; Summary = "Performance enhancements & Video streaming" },
  @{ Row = 33; Date = 45317; Start = 0.020833333333333332; End = 0.14583333333333334; Summary = "Linux bugfixes" },
  @{ Row = 34; Date = 45321; Start = 0.66666666666666663; End = 0.89583333333333337; Summary = "DEM Terrain demo" },
  @{ Row = 35; Date = 45328; Start = 0.20833333333333334; End = 0.5625;               Summary = "Documentation improvements & preparation for user testing" },
  @{ Row = 36; Date = 45332; Start = 0.83333333333333337; End = 0.052083333333333336; Summary = "Implemented image saving" }
)

foreach ($e in $entries) {
  $r = $e.Row
  $ws.Cells.Item($r, 2).Value = $e.Date      # column B - Date
  $ws.Cells.Item($r, 3).Value = $e.Start     # column C - Time Start
  $ws.Cells.Item($r, 4).Value = $e.End       # column D - Time End
  $ws.Cells.Item($r, 6).Value = $e.Summary   # column F - Summary
}

# Reflect where the author was last working in the sheet when they saved.
$ws.Range("F37").Select()
